# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap ranking order of Japon / Austria (Japon's case count overtook Austria) ---
# Row 38 used to be Austria, row 39 used to be Japon.
# Japon moves up to row 38 with refreshed stats; Austria drops to row 39 keeping
# its previous (unchanged) stats.
$ws.Cells.Item(38, 1).Value = "Japon"
$ws.Cells.Item(38, 2).Value = 16120
$ws.Cells.Item(38, 3).Value = 71
$ws.Cells.Item(38, 4).Value = 9868
$ws.Cells.Item(38, 5).Value = 5555
$ws.Cells.Item(38, 6).Value = 259
$ws.Cells.Item(38, 7).Value = 19
$ws.Cells.Item(38, 8).Value = 697

$ws.Cells.Item(39, 1).Value = "Austria"
$ws.Cells.Item(39, 2).Value = 16058
$ws.Cells.Item(39, 3).Value = 61
$ws.Cells.Item(39, 4).Value = 14405
$ws.Cells.Item(39, 5).Value = 1027
$ws.Cells.Item(39, 6).Value = 54
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 626

# --- Refresh case counts for other countries ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1456306
$ws.Cells.Item(4, 3).Value = 25959
$ws.Cells.Item(4, 4).Value = 317857
$ws.Cells.Item(4, 5).Value = 1051554
$ws.Cells.Item(4, 7).Value = 1698
$ws.Cells.Item(4, 8).Value = 86895

# Row 101: Tunez
$ws.Cells.Item(101, 4).Value = 770
$ws.Cells.Item(101, 5).Value = 217
$ws.Cells.Item(101, 6).Value = 3

# Row 125: Jamaica
$ws.Cells.Item(125, 4).Value = 118
$ws.Cells.Item(125, 5).Value = 382

# Row 148: Sudan del Sur
$ws.Cells.Item(148, 4).Value = 3
$ws.Cells.Item(148, 5).Value = 200

# --- Update the "last updated" timestamp string ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Mayo de 2020 a las 02:05"
